# Add a new "2022-Q4" quarterly sheet (cloned from "2022-Q3" so that the
# fund list / styling / layout matches the existing quarterly sheets),
# positioned right after "总计", then update its figures, then update the
# "总计" (totals) sheet with a new leading row for 2022-Q4, shifting the
# rest down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Clone "2022-Q3" -> new sheet placed right after "总计", rename it.
# ---------------------------------------------------------------------
$template = $wb.Sheets("2022-Q3")
$template.Copy($null, $wb.Sheets("总计"))
$q4 = $wb.Sheets("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Helper: write a value as TEXT (matches the source workbook, where these
# numeric-looking figures are stored as strings, not numbers) without
# leaving a stray "quote prefix" cell style behind.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 2) Update the 2022-Q4 figures (fund code/name columns B/C stay as in
#    the template; D/E/F/G are text, H is numeric).
# ---------------------------------------------------------------------
Set-TextValue $q4.Range("D2") "17.90"
Set-TextValue $q4.Range("E2") "94.16"
Set-TextValue $q4.Range("F2") "3.66"
Set-TextValue $q4.Range("G2") "0.6551"
$q4.Range("H2").Value = 5

Set-TextValue $q4.Range("D3") "8.24"
Set-TextValue $q4.Range("E3") "94.16"
Set-TextValue $q4.Range("F3") "3.66"
Set-TextValue $q4.Range("G3") "0.3016"
$q4.Range("H3").Value = 5

Set-TextValue $q4.Range("D4") "6.96"
Set-TextValue $q4.Range("E4") "99.52"
Set-TextValue $q4.Range("F4") "2.67"
Set-TextValue $q4.Range("G4") "0.1858"
$q4.Range("H4").Value = 3

Set-TextValue $q4.Range("D5") "3.50"
Set-TextValue $q4.Range("E5") "98.83"
Set-TextValue $q4.Range("F5") "2.65"
Set-TextValue $q4.Range("G5") "0.0928"
$q4.Range("H5").Value = 3

Set-TextValue $q4.Range("D6") "0.97"
Set-TextValue $q4.Range("E6") "98.95"
Set-TextValue $q4.Range("F6") "2.66"
Set-TextValue $q4.Range("G6") "0.0258"
$q4.Range("H6").Value = 3

# ---------------------------------------------------------------------
# 3) Update "总计" sheet: insert a new row for 2022-Q4 right under the
#    header, pushing the existing rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Sheets("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.26

# Re-number the index column (A) for the rows that shifted down so it
# keeps counting 0,1,2,...
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
$total.Range("A10").Value = 8

# ---------------------------------------------------------------------
# 4) Keep the previously-active tab ("2020-Q4") selected, matching the
#    original workbook's selected tab.
# ---------------------------------------------------------------------
$wb.Sheets("2020-Q4").Activate()
